$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated per-line active-power results (pl_mw) for the 380 kV case.
# Data rows 2-25 correspond to index 0-23 in column A; columns B-F, M, O carry
# the recomputed values while G-L, N remain zero (unchanged).
$cols = @("B","C","D","E","F","M","O")
$data = @(
    @(0.8189485623979635, 0.206113156611309, 0.01661023264770733, 0.1162816279271155, 0.5193643288226966, 0.3402561718116885, 1.682421352132451),
    @(0.7180603288547331, 0.1833884566699169, 0.01506045876318751, 0.1113171366859405, 0.5154294819497309, 0.3024497049344745, 1.683880131782217),
    @(0.6559750088288752, 0.1693662821470525, 0.01410483301199861, 0.1083944561522117, 0.5135291143290246, 0.279281046064753, 1.686508227068231),
    @(0.6306408651642528, 0.1636350846413848, 0.01371441688246477, 0.1072346831979445, 0.5128838973984386, 0.2698508697148299, 1.68801315840787),
    @(0.6264321403134261, 0.1626824039169037, 0.01364952961489507, 0.1070439813154351, 0.5127845492537446, 0.2682856754249983, 1.688289212488456),
    @(0.6556334792271912, 0.1692890577926391, 0.01409957169635589, 0.1083786889747103, 0.5135198902073483, 0.2791538218791487, 1.686526768222791),
    @(0.7841919542180449, 0.1982922165578316, 0.01607673009400656, 0.114543619128014, 0.5179002895121485, 0.3272112240548211, 1.68256377399527),
    @(1.03514700229789, 0.2546071468800619, 0.01992062790949944, 0.1276444601686748, 0.5306039136335627, 0.4218107545084138, 1.688615801237745),
    @(1.218788417262431, 0.2956286025206509, 0.02272324652925306, 0.1379085384689276, 0.5424777852966969, 0.491545864713558, 1.701600527165255),
    @(1.302165886917123, 0.3142115150776874, 0.02399333549279703, 0.1427215359027727, 0.5484382938233807, 0.5233242407134497, 1.709385711006405),
    @(1.333714577942033, 0.3212369096491727, 0.02447356517692612, 0.1445651473780529, 0.5507763043676448, 0.5353660005513348, 1.712605735224287),
    @(1.326921114491995, 0.3197243843482624, 0.024370171858763, 0.1441671515381486, 0.5502691660462062, 0.5327722373072277, 1.711900121390528),
    @(1.304761915840174, 0.3147897319561821, 0.02403285898713392, 0.1428727873517204, 0.5486290191895478, 0.5243147636561361, 1.709645163662174),
    @(1.291185532923805, 0.3117656005384504, 0.02382614975771702, 0.1420827016407529, 0.5476349318981875, 0.5191353557944751, 1.708299406551873),
    @(1.213336136798944, 0.2944125675476812, 0.02264014333565711, 0.1375969218825475, 0.5420995445361001, 0.4894701897874825, 1.701129701764216),
    @(1.165535675899321, 0.2837468334799667, 0.02191130663219099, 0.1348821270269269, 0.5388473054653531, 0.4712857867037314, 1.697213723266032),
    @(1.138026947618641, 0.2776048680865983, 0.02149164558431949, 0.1333341643496269, 0.5370292913401471, 0.4608318252870873, 1.695138047170417),
    @(1.128710398603971, 0.2755240578840414, 0.02134947860522374, 0.1328123613503678, 0.5364227594858306, 0.4572931914514271, 1.694465549761162),
    @(1.170625695537751, 0.2848829788791249, 0.02198893969244864, 0.1351697202160054, 0.5391880655414241, 0.4732210066174929, 1.697612285663325),
    @(1.311271285801297, 0.3162394751710451, 0.02413195587918437, 0.1432523998164399, 0.5491085705602785, 0.5267987121660269, 1.710300104299307),
    @(1.403047656824413, 0.3366652746537682, 0.02552829881911123, 0.1486576641926973, 0.5560639487228087, 0.5618613458505592, 1.720178275538245),
    @(1.354078436813381, 0.3257699272813568, 0.02478344224559237, 0.1457614227133703, 0.5523083993855238, 0.5431435036798717, 1.714760389128429),
    @(1.168324583561343, 0.284369358959907, 0.02195384378720888, 0.1350396594651642, 0.5390338467764622, 0.4723460921623968, 1.697431548524634),
    @(0.9673836328364587, 0.2394336589642592, 0.01888444193939165, 0.1239897549794406, 0.526723462292999, 0.3961792111009004, 1.685487169080346)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $addr = "$($cols[$j])$r"
        $ws.Range($addr).Value = $data[$i][$j]
    }
}
